# Commit: "Fixed typo in Tiers table"
#   "refraction-limited" -> "diffraction-limited" in the Tier-2
#   ("Advanced Quantification and/or Live Cell Imaging") description cell,
#   which is repeated (cell D4) on all three sheets of the workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Tier system_v02-00
$ws2 = $wb.Worksheets.Item(2)   # Tier system_v02-00 SUMMARY
$ws3 = $wb.Worksheets.Item(3)   # Tier system_v02-00 MINIMAL

$newText = "Identification and localization of diffraction-limited particles, super-resolution microscopy, tracking of intracellular dynamics"

$ws1.Range("D4").Value = $newText
$ws2.Range("D4").Value = $newText
$ws3.Range("D4").Value = $newText

# Reflect the cell selections/active sheet left behind by the editing session:
# the SUMMARY sheet's remembered selection moves off D4 onto G4 ...
$ws2.Range("G4").Select() | Out-Null

# ... and the MINIMAL sheet becomes the active tab, with E8 selected there.
$ws3.Activate() | Out-Null
$ws3.Range("E8").Select() | Out-Null
